$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 131
$ws1.Range("H2").Value = 1.79
$ws1.Range("L2").Value = 1.09

$ws1.Range("D3").Value = 127
$ws1.Range("H3").Value = 0.82
$ws1.Range("L3").Value = 1.08

$ws1.Range("L4").Value = 0.93

$ws1.Range("L5").Value = 0.84

$ws1.Range("L6").Value = 0.95

$ws1.Range("L7").Value = 0.92

$ws1.Range("L8").Value = 1.13

$ws1.Range("L9").Value = 1.19

$ws1.Range("L10").Value = 1.12

$ws1.Range("L11").Value = 1.18

$ws1.Range("L12").Value = 0.95

$ws1.Range("L13").Value = 1.03

$ws1.Range("L14").Value = 0.95

$ws1.Range("L15").Value = 0.9399999999999999

$ws1.Range("D16").Value = 105
$ws1.Range("L16").Value = 1.2

$ws1.Range("D17").Value = 97
$ws1.Range("L17").Value = 0.85

# --- Sheet: Summary ---
# These cells hold numeric-looking values that are stored as TEXT
# (inlineStr) in the source file, so a leading apostrophe is used to force
# text entry and keep the cell type as text rather than becoming numeric.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'1790"
$ws2.Range("B10").Value = "'950"
$ws2.Range("B11").Value = "'488"
$ws2.Range("B12").Value = "'131"
$ws2.Range("B14").Value = "'98"
